$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "with" " " "one" " " "Edit" " " "Text" " " "so"  ->  merge into a single
#    run "with one Edit Text so" (same formatting throughout, so a normal
#    Find/Replace over the whole phrase collapses the adjacent runs into one).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("with one Edit Text so", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "with one Edit Text so", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "a" " " "button"  ->  merge into a single run "a button".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("a button", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a button", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "com.keshavrkaranth.text_to_speech" -> "com.vinayak.text_to_speech",
#    split across three runs: "com.", "vinayak", ".text_to_speech" (the
#    middle run is the only part whose text actually changed).
# ---------------------------------------------------------------------------

# First swap the package-name segment; this temporarily collapses the run
# into "com.vinayak.text_to_speech" as one run.
$d.Content.Find.Execute("keshavrkaranth", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "vinayak", 2) | Out-Null

# Now force Word to split that run around "vinayak" by nudging a character
# property on just that sub-range (bold on, then back off) - this leaves the
# run boundaries in place (com. | vinayak | .text_to_speech) without changing
# the visible formatting.
$rngSplit = $d.Content
$rngSplit.Find.Execute("vinayak", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$rngSplit.Font.Bold = 1
$rngSplit.Font.Bold = 0

Write-Output "done"
